$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '42.006.21'
$ws.Range("E2").Value2 = '  -1.34%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '2.299.87'
$ws.Range("E3").Value2 = '  -2.99%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.00'
$ws.Range("E4").Value2 = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '310.77'
$ws.Range("E5").Value2 = '  -6.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '105.10'
$ws.Range("E6").Value2 = '  +4.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.624'
$ws.Range("E7").Value2 = '  -2.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '1.00'
$ws.Range("E8").Value2 = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.605'
$ws.Range("E9").Value2 = '  -4.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '39.74'
$ws.Range("E10").Value2 = '  -1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0909'
$ws.Range("E11").Value2 = '  -1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '8.26'
$ws.Range("E12").Value2 = '  -2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.105'
$ws.Range("E13").Value2 = '  -0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '0.964'
$ws.Range("E14").Value2 = '  -4.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '15.37'
$ws.Range("E15").Value2 = '  -6.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '2.655.29'
$ws.Range("E16").Value2 = '  -2.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '2.318.67'
$ws.Range("E17").Value2 = '  -2.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '42.030.69'
$ws.Range("E18").Value2 = '  -1.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '7.56'
$ws.Range("E19").Value2 = '  -4.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '0.0000104'
$ws.Range("E20").Value2 = '  -2.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '74.13'
$ws.Range("E21").Value2 = '  -1.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '3.46'
$ws.Range("E22").Value2 = '  -9.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '258.52'
$ws.Range("E23").Value2 = '  -5.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.26'
$ws.Range("E24").Value2 = '  -2.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '9.16'
$ws.Range("E25").Value2 = '  -6.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '1.01'
$ws.Range("E26").Value2 = '  +0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '10.95'
$ws.Range("E27").Value2 = '  -4.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '2.26'
$ws.Range("E28").Value2 = '  +2.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '22.71'
$ws.Range("E29").Value2 = '  -2.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '35.50'
$ws.Range("E30").Value2 = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '164.16'
$ws.Range("E31").Value2 = '  -6.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '0.0883'
$ws.Range("E32").Value2 = '  -2.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '2.90'
$ws.Range("E33").Value2 = '  -6.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '5.82'
$ws.Range("E34").Value2 = '  -4.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '0.128'
$ws.Range("E35").Value2 = '  -4.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.116'
$ws.Range("E36").Value2 = '  +9.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '4.49'
$ws.Range("E37").Value2 = '  -3.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.0349'
$ws.Range("E38").Value2 = '  -3.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.64'
$ws.Range("E39").Value2 = '  -6.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '2.69'
$ws.Range("E40").Value2 = '  -7.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '98.24'
$ws.Range("E41").Value2 = '  +9.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '1.46'
$ws.Range("E42").Value2 = '  -4.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '69.64'
$ws.Range("E43").Value2 = '  -1.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.229'
$ws.Range("E44").Value2 = '  -2.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '1.00'
$ws.Range("E45").Value2 = '  +0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '12.09'
$ws.Range("E46").Value2 = '  -0.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '110.63'
$ws.Range("E47").Value2 = '  -6.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '5.34'
$ws.Range("E48").Value2 = '  -2.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '8.94'
$ws.Range("E49").Value2 = '  -2.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '73.00'
$ws.Range("E50").Value2 = '  +3.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '1.25'
$ws.Range("E51").Value2 = '  -1.57%  '
